$wb = $excel.ActiveWorkbook

$oldDate = "Date:                Thu, 02 Jan 2020"
$newDate = "Date:                Sun, 05 Jan 2020"
$oldTime = "Time:                        20:48:41"
$newTime = "Time:                        21:22:19"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value2
    if ($text -ne $null -and $text -like "*Date:*") {
        $text = $text.Replace($oldDate, $newDate)
        $text = $text.Replace($oldTime, $newTime)
        $cell.Value2 = $text
    }
}
